$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New LEM/LED pairs to append below the existing table (rows 38-48)
$newRows = @(
  @("LEM-198-32-3022KS", "LED-198-S70-3022"),
  @("LEM-219-00-6022KS", "LED-219-S00-6022"),
  @("LEM-236-00-35KS", "LED-236-S00-35"),
  @("LEM-239-00-35KS", "LED-239-S00-35"),
  @("LEM-240-01-35KH", "LED-240-H01-35"),
  @("LEM-274-00-27KH", "LED-274-H00-27"),
  @("LEM-274-00-30KH", "LED-274-H00-30"),
  @("LEM-275-32-2722KS", "LED-275-S00-2722"),
  @("LEM-275-32-3522KS", "LED-275-S00-3522"),
  @("LEM-293-00-30KH", "LED-293-H00-30"),
  @("LEM-313-00-3022KH", "LED-313-H00-3022")
)

$startRow = 38
for ($i = 0; $i -lt $newRows.Count; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
  $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Re-sort the data range (excluding header row) by column A ascending,
# using the worksheet Sort object so the sortState is persisted like Excel's
# Data > Sort dialog would record it.
$sortRange = $ws.Range("A1:B48")
$keyRange = $ws.Range("A1:A48")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, 0, 1, $null, 0)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# Update the selected cell as recorded in the saved workbook view
$ws.Range("D3").Select()
